# Generate Report for Handoff
#
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps for the row corresponding to
# d81a3c54-a465-46b1-9302-12fe23154d8d (the file that was just handed off)
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date", row 7
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-15 16:40:26"

# zh-cn sheet: column H = "Latest Handoff Datetime", row 7
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-15 16:40:21"

# de-de sheet: column H = "Latest Handoff Datetime", row 7
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-15 16:40:26"
